$d = $word.ActiveDocument

# wdReplaceAll = 2, wdFindContinue = 1 (wrap within the search scope)
$wdFindContinue = 1
$wdReplaceAll = 2

# Helper: find the paragraph whose text contains a given marker string and
# return its Range (duplicated so later Find calls don't walk past it).
function Get-ParagraphRangeContaining($doc, [string]$marker) {
    $paras = $doc.Paragraphs
    for ($i = 1; $i -le $paras.Count; $i++) {
        $p = $paras.Item($i)
        if ($p.Range.Text -like "*$marker*") {
            return $p.Range.Duplicate
        }
    }
    return $null
}

# --- Edit 1: Professional summary paragraph ---
# "...errors affecting all Black and Asian-American voters, developed..."
# -> "...errors affecting 50M voters, developed..."
$summaryRange = Get-ParagraphRangeContaining $d "Data engineering professional with 15+ years"
$summaryRange.Find.Execute(
    "affecting all Black and Asian-American voters,", $true, $false, $false,
    $false, $false, $true, $wdFindContinue, $false,
    "affecting 50M voters,", $wdReplaceAll)

# --- Edit 2: Bullet point under Siege Analytics experience ---
# "...errors affecting all Black and Asian-American voters, developed geospatial machine learning..."
# -> "...errors affecting " + bold/colored "50M" + " voters, developed geospatial machine learning..."
$bulletRange = Get-ParagraphRangeContaining $d "Discovered systematic race coding errors affecting"
$bulletRange.Find.Execute(
    "affecting all Black and Asian-American voters,", $true, $false, $false,
    $false, $false, $true, $wdFindContinue, $false,
    "affecting 50M voters,", $wdReplaceAll)

# Re-acquire the (now-updated) paragraph range and bold/color just the "50M" token.
$bulletRange2 = Get-ParagraphRangeContaining $d "Discovered systematic race coding errors affecting 50M voters,"
$bulletRange2.Find.Execute("50M")
$bulletRange2.Font.Bold = 1
$bulletRange2.Font.Color = 5258796

# --- Edit 3: Impact line under the project section ---
# "Impact: Corrected demographic data affecting all Black and Asian-American voters, improved..."
# -> "Impact: Corrected demographic data affecting 50M voters nationwide, improved..."
$impactRange = Get-ParagraphRangeContaining $d "Impact: Corrected demographic data affecting"
$impactRange.Find.Execute(
    "affecting all Black and Asian-American voters, improved",
    $true, $false, $false, $false, $false, $true, $wdFindContinue, $false,
    "affecting 50M voters nationwide, improved", $wdReplaceAll)
